$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 147398.9
$ws.Range("J17").Value = 150726.16
$ws.Range("L17").Value = 452178.48
$ws.Range("N17").Value = -452514.48

# Row 20
$ws.Range("H20").Value = 1000
$ws.Range("J20").Value = 1000
$ws.Range("L20").Value = 1000
$ws.Range("N20").Value = -1460

# Row 29
$ws.Range("H29").Value = 700
$ws.Range("I29").Value = 700
$ws.Range("K29").Value = 2100
$ws.Range("M29").Value = -1819

# Row 35
$ws.Range("H35").Value = 1000
$ws.Range("J35").Value = 1000
$ws.Range("L35").Value = 1000
$ws.Range("N35").Value = -1758

# Row 38
$ws.Range("H38").Value = 733.5
$ws.Range("I38").Value = 119.85714
$ws.Range("K38").Value = 359.57142
$ws.Range("M38").Value = 12.42858000000001

# Row 55
$ws.Range("H55").Value = 1027
$ws.Range("J55").Value = 2115.6667
$ws.Range("L55").Value = 2115.6667
$ws.Range("N55").Value = -2543.6667

# Row 62
$ws.Range("H62").Value = 2407
$ws.Range("J62").Value = 3166.3333
$ws.Range("L62").Value = 3166.3333
$ws.Range("N62").Value = -4414.3333

# Row 65
$ws.Range("H65").Value = 2407
$ws.Range("J65").Value = 3166.3333
$ws.Range("L65").Value = 15831.6665
$ws.Range("N65").Value = -22071.6665

# Row 80
$ws.Range("H80").Value = 1297.5
$ws.Range("J80").Value = 1647.5
$ws.Range("L80").Value = 4942.5
$ws.Range("N80").Value = -6938.5

# Row 83
$ws.Range("H83").Value = 1297.5
$ws.Range("J83").Value = 1647.5
$ws.Range("L83").Value = 14827.5
$ws.Range("N83").Value = -24811.5

# Row 111
$ws.Range("H111").Value = 910.25
$ws.Range("I111").Value = 905.8333
$ws.Range("J111").Value = 923.5
$ws.Range("K111").Value = 2717.4999
$ws.Range("L111").Value = 2770.5
$ws.Range("M111").Value = 349.5001000000002
$ws.Range("N111").Value = -8904.5

# Row 116
$ws.Range("H116").Value = 5900
$ws.Range("I116").Value = 5900
$ws.Range("K116").Value = 5900
$ws.Range("M116").Value = -2458

# Row 133
$ws.Range("H133").Value = 141654.17
$ws.Range("J133").Value = 141654.17
$ws.Range("L133").Value = 141654.17
$ws.Range("N133").Value = -151774.17

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 926.625
$ws.Range("J4").Value = 796.3333
$ws.Range("L4").Value = 796.3333
$ws.Range("N4").Value = -1028.3333

# Row 35
$ws.Range("H35").Value = 5927.2
$ws.Range("I35").Value = 5984
$ws.Range("J35").Value = 5700
$ws.Range("K35").Value = 5984
$ws.Range("L35").Value = 5700
$ws.Range("M35").Value = -5578
$ws.Range("N35").Value = -6512

# Row 132
$ws.Range("H132").Value = 2086.7964
$ws.Range("I132").Value = 1706.9348
$ws.Range("J132").Value = 4271
$ws.Range("K132").Value = 5120.8044
$ws.Range("L132").Value = 12813
$ws.Range("M132").Value = -2590.8044
$ws.Range("N132").Value = -17873

$ws = $wb.Worksheets.Item("BSM")
# Row 35
$ws.Range("H35").Value = 60000
$ws.Range("J35").Value = 60000
$ws.Range("L35").Value = 60000
$ws.Range("N35").Value = -60620

# Row 99
$ws.Range("H99").Value = 6220.2
$ws.Range("I99").Value = 6908.857
$ws.Range("K99").Value = 6908.857
$ws.Range("M99").Value = -5410.857

# Row 105
$ws.Range("H105").Value = 40499.332
$ws.Range("I105").Value = 99999
$ws.Range("J105").Value = 10749.5
$ws.Range("K105").Value = 99999
$ws.Range("L105").Value = 10749.5
$ws.Range("M105").Value = -98252
$ws.Range("N105").Value = -14243.5

$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 2157.4
$ws.Range("I5").Value = 3401
$ws.Range("J5").Value = 292
$ws.Range("K5").Value = 3401
$ws.Range("L5").Value = 292
$ws.Range("M5").Value = -3289
$ws.Range("N5").Value = -516

# Row 22
$ws.Range("H22").Value = 467.0926
$ws.Range("I22").Value = 448
$ws.Range("J22").Value = 791.6667
$ws.Range("K22").Value = 448
$ws.Range("L22").Value = 791.6667
$ws.Range("M22").Value = -98
$ws.Range("N22").Value = -1491.6667

# Row 122
$ws.Range("H122").Value = 2126.3125
$ws.Range("I122").Value = 2128.5386
$ws.Range("K122").Value = 6385.6158
$ws.Range("M122").Value = -3935.6158

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 5633.222
$ws.Range("I32").Value = 2924.75
$ws.Range("J32").Value = 7800
$ws.Range("K32").Value = 8774.25
$ws.Range("L32").Value = 23400
$ws.Range("M32").Value = -8491.25
$ws.Range("N32").Value = -23966

# Row 129
$ws.Range("H129").Value = 10233.529
$ws.Range("I129").Value = 893.9
$ws.Range("J129").Value = 23575.857
$ws.Range("K129").Value = 2681.7
$ws.Range("L129").Value = 70727.571
$ws.Range("M129").Value = 2318.3
$ws.Range("N129").Value = -80727.571

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 82.916664
$ws.Range("I2").Value = 73
$ws.Range("K2").Value = 73
$ws.Range("M2").Value = 40

# Row 21
$ws.Range("H21").Value = 2513875
$ws.Range("I21").Value = 5008000
$ws.Range("J21").Value = 19750
$ws.Range("K21").Value = 5008000
$ws.Range("L21").Value = 19750
$ws.Range("M21").Value = -5007827
$ws.Range("N21").Value = -20096

# Row 30
$ws.Range("H30").Value = 2513875
$ws.Range("I30").Value = 5008000
$ws.Range("J30").Value = 19750
$ws.Range("K30").Value = 5008000
$ws.Range("L30").Value = 19750
$ws.Range("M30").Value = -5007895
$ws.Range("N30").Value = -19960

# Row 63
$ws.Range("H63").Value = 61111
$ws.Range("I63").Value = 100000
$ws.Range("J63").Value = 22222
$ws.Range("K63").Value = 100000
$ws.Range("L63").Value = 22222
$ws.Range("M63").Value = -99314
$ws.Range("N63").Value = -23594

# Row 66
$ws.Range("H66").Value = 61111
$ws.Range("I66").Value = 100000
$ws.Range("J66").Value = 22222
$ws.Range("K66").Value = 300000
$ws.Range("L66").Value = 66666
$ws.Range("M66").Value = -296568
$ws.Range("N66").Value = -73530

# Row 126
$ws.Range("H126").Value = 4849.7
$ws.Range("I126").Value = 4624.5
$ws.Range("J126").Value = 4999.8335
$ws.Range("K126").Value = 13873.5
$ws.Range("L126").Value = 14999.5005
$ws.Range("M126").Value = -11403.5
$ws.Range("N126").Value = -19939.5005

$ws = $wb.Worksheets.Item("LTW")
# Row 14
$ws.Range("H14").Value = 82175
$ws.Range("I14").Value = 82175
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 82175
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -82003
$ws.Range("N14").ClearContents()

# Row 55
$ws.Range("H55").Value = 794.1
$ws.Range("I55").Value = 518.45
$ws.Range("K55").Value = 518.45
$ws.Range("M55").Value = -345.45

# Row 61
$ws.Range("H61").Value = 2906.1304
$ws.Range("I61").Value = 2400.389
$ws.Range("J61").Value = 4726.8
$ws.Range("K61").Value = 2400.389
$ws.Range("L61").Value = 4726.8
$ws.Range("M61").Value = -2198.389
$ws.Range("N61").Value = -5130.8

# Row 113
$ws.Range("H113").Value = 2906.1304
$ws.Range("I113").Value = 2400.389
$ws.Range("J113").Value = 4726.8
$ws.Range("K113").Value = 2400.389
$ws.Range("L113").Value = 4726.8
$ws.Range("M113").Value = -230.3890000000001
$ws.Range("N113").Value = -9066.799999999999

# Row 132
$ws.Range("H132").Value = 2949.0557
$ws.Range("J132").Value = 3185.3333
$ws.Range("L132").Value = 9555.999899999999
$ws.Range("N132").Value = -14615.9999

# Row 136
$ws.Range("H136").Value = 13309.759
$ws.Range("I136").Value = 43511.582
$ws.Range("K136").Value = 130534.746
$ws.Range("M136").Value = -127984.746

$ws = $wb.Worksheets.Item("WVR")
# Row 10
$ws.Range("H10").Value = 22496.5
$ws.Range("J10").Value = 22496.5
$ws.Range("L10").Value = 22496.5
$ws.Range("N10").Value = -22834.5

# Row 34
$ws.Range("H34").Value = 19999.666
$ws.Range("I34").Value = 19999.666
$ws.Range("K34").Value = 19999.666
$ws.Range("M34").Value = -19796.666

# Row 96
$ws.Range("H96").Value = 27417.375
$ws.Range("I96").Value = 4098
$ws.Range("K96").Value = 4098
$ws.Range("M96").Value = -2725

# Row 107
$ws.Range("H107").Value = 454.21054
$ws.Range("I107").Value = 470.46155
$ws.Range("J107").Value = 419
$ws.Range("K107").Value = 1411.38465
$ws.Range("L107").Value = 1257
$ws.Range("M107").Value = 508.61535
$ws.Range("N107").Value = -5097

# Row 113
$ws.Range("H113").Value = 408.76923
$ws.Range("I113").Value = 393.3889
$ws.Range("K113").Value = 1180.1667
$ws.Range("M113").Value = 989.8333
